$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note under the book summary table
$ws.Range("B8").Value = "(not the pdf pages)"

# Book 3 chapter page numbers (column J), rows 12-150
$bookThreePages = @(
    14,
    27,
    35,
    47,
    61,
    79,
    87,
    96,
    113,
    125,
    131,
    147,
    166,
    187,
    200,
    205,
    217,
    239,
    250,
    273,
    287,
    292,
    305,
    313,
    319,
    334,
    347,
    365,
    378,
    393,
    415,
    419,
    436,
    442,
    445,
    452,
    461,
    468,
    484,
    498,
    505,
    530,
    550,
    563,
    575,
    582,
    590,
    599,
    614,
    623,
    638,
    645,
    652,
    667,
    680,
    689,
    699,
    712,
    722,
    732,
    740,
    749,
    760,
    771,
    781,
    789,
    807,
    817,
    833,
    847,
    858,
    869,
    886,
    897,
    912,
    924,
    934,
    943,
    955,
    969,
    982,
    992,
    1002,
    1015,
    1030,
    1043,
    1051,
    1057,
    1061,
    1073,
    1091,
    1106,
    1110,
    1116,
    1121,
    1129,
    1133,
    1135,
    1139,
    1145,
    1154,
    1167,
    1174,
    1181,
    1193,
    1201,
    1205,
    1211,
    1226,
    1241,
    1251,
    1264,
    1277,
    1286,
    1302,
    1313,
    1323,
    1333,
    1340,
    1362,
    1382,
    1396,
    1403,
    1418,
    1423,
    1431,
    1434,
    1448,
    1452,
    1464,
    1480,
    1499,
    1519,
    1533,
    1547,
    1608,
    1637,
    1662,
    1668
)

for ($i = 0; $i -lt $bookThreePages.Length; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 10).Value = $bookThreePages[$i]
}

# Restore sheet view state (scroll position / selection) as recorded at save time
$ws.Range("J158").Select() | Out-Null
